# "modificar desplegable de servicios valorados"
# Give each valued service (gestoria) row its own distinct NIF instead of
# reusing the same placeholder value, flip row 3's "Activa" flag to FALSE,
# and widen the NIF/Activa columns so the new values are readable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Distinct NIF per row (was "00000000T" for every row).
$ws.Range("N3").Value = "00000002T"
$ws.Range("N4").Value = "00000003T"
$ws.Range("N5").Value = "00000004T"

# Row 3 is no longer an active/valued service.
$ws.Range("M3").Value = $false

# Widen the "Activa" and "NIF" columns so the dropdown/values fit.
$ws.Columns.Item(13).ColumnWidth = 14
$ws.Columns.Item(14).ColumnWidth = 23

# Leave the cursor where the edit was made.
$ws.Range("N8").Select() | Out-Null
